# "support selling in benchmark"
# Edit the test_transactions.xlsx benchmark fixture:
#  - Row 29 becomes a new Cash withdrawal transaction (company_a / Cash / -1000)
#    instead of the old company_b SELL NFLX row, and gets a note "withdrawal".
#  - G24:G28 formulas become one shared-formula group (re-entering the identical
#    formula across the whole range at once, same as Excel's own fill-down,
#    makes Excel store it as a single shared-formula group).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$quote = [char]34

# Re-apply the (unchanged) G24:G28 cost formula across the whole block in one
# shot so Excel collapses it into a single shared formula (si group), matching
# how Excel re-writes a column of identical relative formulas.
$ws.Range("G24:G28").Formula = "=IF(C24=" + $quote + "Cash" + $quote + ",1,-1)*E24*F24"

# Update row 29 in place to the new Cash withdrawal transaction.
$ws.Range("A29").Value = "company_a"
$ws.Range("B29").Value = 44704
$ws.Range("C29").Value = "Cash"
$ws.Range("D29").Value = "Cash"
$ws.Range("E29").Value = -1000
$ws.Range("F29").Value = 1
$ws.Range("G29").Formula = "=IF(C29=" + $quote + "Cash" + $quote + ",1,-1)*E29*F29"
$ws.Range("H29").Value = "withdrawal"

# Move the visible selection down to the newly-edited row, like the author
# scrolling to/clicking on the new transaction before saving.
[void]$ws.Range("L36").Select()
